# Allineamento ambiente di test con daf+mur
# Insert a new licence row (B.1.17 - "Dichiarazioni di uso standard beni
# culturali (BCS)") right after the other "B.1" rows, before the "B.2" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 currently holds the first "B.2" record; push it (and everything
# below) down by one row to make room for the new "B.1.17" record.
$ws.Rows.Item(42).Insert()

$ws.Cells.Item(42, 1).Value = "B"
$ws.Cells.Item(42, 2).Value = "Licenza Non Aperta"
$ws.Cells.Item(42, 3).Value = "B.1"
$ws.Cells.Item(42, 4).Value = "Solo uso non commerciale"
$ws.Cells.Item(42, 5).Value = "B.1.17"
$ws.Cells.Item(42, 6).Value = "Dichiarazioni di uso standard beni culturali (BCS)"

# Match formatting used by the other repeated "B.1" rows: vertically
# centred labels, with column C (the repeated "B.1" code) rendered in the
# plain-black "repeat" font used throughout the B.1 block.
$ws.Cells.Item(42, 2).VerticalAlignment = -4108
$ws.Cells.Item(42, 4).VerticalAlignment = -4108

$ws.Cells.Item(42, 3).HorizontalAlignment = -4131
$ws.Cells.Item(42, 3).VerticalAlignment = -4108
$ws.Cells.Item(42, 3).Font.Color = 0

$ws.Cells.Item(42, 5).HorizontalAlignment = -4131
$ws.Cells.Item(42, 5).VerticalAlignment = -4108

$ws.Range("F41").Select()
